$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)
$ws.Activate()

# Insert a new row above row 142, shifting existing rows 142-205 down to 143-206
$ws.Rows(142).Insert()

# Populate the newly inserted row 142 with the new DRC rule entry
$ws.Range("F142").Value = "MetalxMinSpace11"
$ws.Range("G142").Value = 500
$ws.Range("H142").Value = "GR604e2_M6"
$ws.Range("I142").Value = "M6 minimum space to (M6 with width > 1.500), for run length > 1.500, >= 0.5"

# Update view state to match: active window scrolled and new selection
$excel.ActiveWindow.ScrollRow = 46
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F141").Select()
